# Applies the PSP sheet update: fills in the previously-empty Dec 5/6/7
# log rows (A31:F33) on the "작성자명" (first) worksheet with their real
# date / start-stop time / interruption / delta / activity data, and
# leaves the selection on E33 (matching the author's last-saved cursor
# position in the source diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 31 : 12월 5일 -------------------------------------------------
$ws.Range("A31").Value = "12월 5일"
# "12" keeps the cell's base font; "월 5일" is an explicit Arial Unicode MS run
$ws.Range("A31").Characters(3, 4).Font.Name = "Arial Unicode MS"

$ws.Range("B31").Value = 0.625
$ws.Range("C31").Value = 0.6875
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 90

$ws.Range("F31").Value = "추천로직기획"
$ws.Range("F31").Font.Name = "Arial Unicode MS"

# --- Row 32 : 12월 6일 -------------------------------------------------
$ws.Range("A32").Value = "12월 6일"
$ws.Range("A32").Characters(3, 4).Font.Name = "Arial Unicode MS"

$ws.Range("B32").Value = 0.54166666666666663
$ws.Range("C32").Value = 0.72916666666666663
$ws.Range("D32").Value = 90
$ws.Range("E32").Value = 180

$ws.Range("F32").Value = "initial data 재정리"
$ws.Range("F32").Font.Name = "Arial Unicode MS"

# --- Row 33 : 12월7일 --------------------------------------------------
$ws.Range("A33").Value = "12월7일"
$ws.Range("A33").Characters(3, 3).Font.Name = "Arial Unicode MS"

$ws.Range("B33").Value = 0.66666666666666663
$ws.Range("C33").Value = 0.91666666666666663
$ws.Range("D33").Value = 120
$ws.Range("E33").Value = 240

$ws.Range("F33").Value = "sequence diagram 설계, srs수정, class design 수정"
# "sequence diagram " stays the default (Arial) run; the Korean tail is 돋움
$ws.Range("F33").Characters(18, 26).Font.Name = "돋움"

# --- Selection / view ---------------------------------------------------
$ws.Activate()
$ws.Range("E33").Select()

Write-Output "PSP sheet rows 31-33 populated"
